# Update Ctf1-Il6st LR-pairs sheet with new TPM-derived values.
# The sending/target cluster set now includes "Resolving-Mac" (previously only
# "Inflammatory-Mac" appeared under Target cluster); the data now contains the
# full 4x5 cross join of clusters x the 4 original target clusters plus the new
# "Resolving-Mac" sending-cluster rows (rows 17-21), and many numeric columns
# (G-J, M-T for the Resolving-Mac target + recomputed edge stats) have new TPM values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Ctf1"
$ws.Range("C2").Value = "Il6st"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.055246
$ws.Range("H2").Value = 0.165738
$ws.Range("I2").Value = 0.04876749798589627
$ws.Range("J2").Value = 0.04876749798589627
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 108.845309
$ws.Range("N2").Value = 326.535927
$ws.Range("O2").Value = 0.3930063530400584
$ws.Range("P2").Value = 0.3930063530400583
$ws.Range("Q2").Value = 6.013267941014
$ws.Range("R2").Value = 54.119411469126
$ws.Range("S2").Value = 0.01916593653032549
$ws.Range("T2").Value = 0.01916593653032548

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Ctf1"
$ws.Range("C3").Value = "Il6st"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.055246
$ws.Range("H3").Value = 0.165738
$ws.Range("I3").Value = 0.04876749798589627
$ws.Range("J3").Value = 0.04876749798589627
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 125.002688
$ws.Range("N3").Value = 375.008064
$ws.Range("O3").Value = 0.4513455929560021
$ws.Range("P3").Value = 0.451345592956002
$ws.Range("Q3").Value = 6.905898501247999
$ws.Range("R3").Value = 62.153086511232
$ws.Range("S3").Value = 0.02201099529542499
$ws.Range("T3").Value = 0.02201099529542498

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Ctf1"
$ws.Range("C4").Value = "Il6st"
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.055246
$ws.Range("H4").Value = 0.165738
$ws.Range("I4").Value = 0.04876749798589627
$ws.Range("J4").Value = 0.04876749798589627
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 18.88973866666667
$ws.Range("N4").Value = 56.669216
$ws.Range("O4").Value = 0.06820493571538706
$ws.Range("P4").Value = 0.06820493571538705
$ws.Range("Q4").Value = 1.043582502378667
$ws.Range("R4").Value = 9.392242521407999
$ws.Range("S4").Value = 0.003326184065128324
$ws.Range("T4").Value = 0.003326184065128322

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Ctf1"
$ws.Range("C5").Value = "Il6st"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.055246
$ws.Range("H5").Value = 0.165738
$ws.Range("I5").Value = 0.04876749798589627
$ws.Range("J5").Value = 0.04876749798589627
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 13.37928666666667
$ws.Range("N5").Value = 40.13786
$ws.Range("O5").Value = 0.04830841776694433
$ws.Range("P5").Value = 0.04830841776694433
$ws.Range("Q5").Value = 0.7391520711866667
$ws.Range("R5").Value = 6.652368640680001
$ws.Range("S5").Value = 0.002355880666151294
$ws.Range("T5").Value = 0.002355880666151293

# Row 6
$ws.Range("A6").Value = "ECs"
$ws.Range("B6").Value = "Ctf1"
$ws.Range("C6").Value = "Il6st"
$ws.Range("D6").Value = "Resolving-Mac"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.055246
$ws.Range("H6").Value = 0.165738
$ws.Range("I6").Value = 0.04876749798589627
$ws.Range("J6").Value = 0.04876749798589627
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 10.83857433333333
$ws.Range("N6").Value = 32.515723
$ws.Range("O6").Value = 0.03913470052160829
$ws.Range("P6").Value = 0.03913470052160829
$ws.Range("Q6").Value = 0.5987878776193333
$ws.Range("R6").Value = 5.389090898574
$ws.Range("S6").Value = 0.001908501428866186
$ws.Range("T6").Value = 0.001908501428866186

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Ctf1"
$ws.Range("C7").Value = "Il6st"
$ws.Range("D7").Value = "ECs"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.6260123333333333
$ws.Range("H7").Value = 1.878037
$ws.Range("I7").Value = 0.5526020925493169
$ws.Range("J7").Value = 0.5526020925493168
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 108.845309
$ws.Range("N7").Value = 326.535927
$ws.Range("O7").Value = 0.3930063530400584
$ws.Range("P7").Value = 0.3930063530400583
$ws.Range("Q7").Value = 68.13850585947766
$ws.Range("R7").Value = 613.246552735299
$ws.Range("S7").Value = 0.2171761330751119
$ws.Range("T7").Value = 0.2171761330751118

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Ctf1"
$ws.Range("C8").Value = "Il6st"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.6260123333333333
$ws.Range("H8").Value = 1.878037
$ws.Range("I8").Value = 0.5526020925493169
$ws.Range("J8").Value = 0.5526020925493168
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 125.002688
$ws.Range("N8").Value = 375.008064
$ws.Range("O8").Value = 0.4513455929560021
$ws.Range("P8").Value = 0.451345592956002
$ws.Range("Q8").Value = 78.25322438781866
$ws.Range("R8").Value = 704.279019490368
$ws.Range("S8").Value = 0.249414519130399
$ws.Range("T8").Value = 0.2494145191303989

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Ctf1"
$ws.Range("C9").Value = "Il6st"
$ws.Range("D9").Value = "Inflammatory-Mac"
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.6260123333333333
$ws.Range("H9").Value = 1.878037
$ws.Range("I9").Value = 0.5526020925493169
$ws.Range("J9").Value = 0.5526020925493168
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 18.88973866666667
$ws.Range("N9").Value = 56.669216
$ws.Range("O9").Value = 0.06820493571538706
$ws.Range("P9").Value = 0.06820493571538705
$ws.Range("Q9").Value = 11.82520937877689
$ws.Range("R9").Value = 106.426884408992
$ws.Range("S9").Value = 0.03769019019851453
$ws.Range("T9").Value = 0.03769019019851452

# Row 10
$ws.Range("A10").Value = "FAPs"
$ws.Range("B10").Value = "Ctf1"
$ws.Range("C10").Value = "Il6st"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.6260123333333333
$ws.Range("H10").Value = 1.878037
$ws.Range("I10").Value = 0.5526020925493169
$ws.Range("J10").Value = 0.5526020925493168
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 13.37928666666667
$ws.Range("N10").Value = 40.13786
$ws.Range("O10").Value = 0.04830841776694433
$ws.Range("P10").Value = 0.04830841776694433
$ws.Range("Q10").Value = 8.375598464535555
$ws.Range("R10").Value = 75.38038618082001
$ws.Range("S10").Value = 0.02669533274576004
$ws.Range("T10").Value = 0.02669533274576003

# Row 11
$ws.Range("A11").Value = "FAPs"
$ws.Range("B11").Value = "Ctf1"
$ws.Range("C11").Value = "Il6st"
$ws.Range("D11").Value = "Resolving-Mac"
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 0.6666666666666666
$ws.Range("G11").Value = 0.6260123333333333
$ws.Range("H11").Value = 1.878037
$ws.Range("I11").Value = 0.5526020925493169
$ws.Range("J11").Value = 0.5526020925493168
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 10.83857433333333
$ws.Range("N11").Value = 32.515723
$ws.Range("O11").Value = 0.03913470052160829
$ws.Range("P11").Value = 0.03913470052160829
$ws.Range("Q11").Value = 6.785081208416777
$ws.Range("R11").Value = 61.065730875751
$ws.Range("S11").Value = 0.02162591739953159
$ws.Range("T11").Value = 0.02162591739953158

# Row 12
$ws.Range("A12").Value = "MuSCs"
$ws.Range("B12").Value = "Ctf1"
$ws.Range("C12").Value = "Il6st"
$ws.Range("D12").Value = "ECs"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.4378713333333333
$ws.Range("H12").Value = 1.313614
$ws.Range("I12").Value = 0.3865237187563814
$ws.Range("J12").Value = 0.3865237187563814
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 108.845309
$ws.Range("N12").Value = 326.535927
$ws.Range("O12").Value = 0.3930063530400584
$ws.Range("P12").Value = 0.3930063530400583
$ws.Range("Q12").Value = 47.66024057890866
$ws.Range("R12").Value = 428.942165210178
$ws.Range("S12").Value = 0.1519062770719267
$ws.Range("T12").Value = 0.1519062770719266

# Row 13
$ws.Range("A13").Value = "MuSCs"
$ws.Range("B13").Value = "Ctf1"
$ws.Range("C13").Value = "Il6st"
$ws.Range("D13").Value = "FAPs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.4378713333333333
$ws.Range("H13").Value = 1.313614
$ws.Range("I13").Value = 0.3865237187563814
$ws.Range("J13").Value = 0.3865237187563814
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 125.002688
$ws.Range("N13").Value = 375.008064
$ws.Range("O13").Value = 0.4513455929560021
$ws.Range("P13").Value = 0.451345592956002
$ws.Range("Q13").Value = 54.73509366481066
$ws.Range("R13").Value = 492.6158429832959
$ws.Range("S13").Value = 0.1744557770336579
$ws.Range("T13").Value = 0.1744557770336579

# Row 14
$ws.Range("A14").Value = "MuSCs"
$ws.Range("B14").Value = "Ctf1"
$ws.Range("C14").Value = "Il6st"
$ws.Range("D14").Value = "Inflammatory-Mac"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.4378713333333333
$ws.Range("H14").Value = 1.313614
$ws.Range("I14").Value = 0.3865237187563814
$ws.Range("J14").Value = 0.3865237187563814
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 18.88973866666667
$ws.Range("N14").Value = 56.669216
$ws.Range("O14").Value = 0.06820493571538706
$ws.Range("P14").Value = 0.06820493571538705
$ws.Range("Q14").Value = 8.271275056291554
$ws.Range("R14").Value = 74.44147550662399
$ws.Range("S14").Value = 0.02636282539025134
$ws.Range("T14").Value = 0.02636282539025133

# Row 15
$ws.Range("A15").Value = "MuSCs"
$ws.Range("B15").Value = "Ctf1"
$ws.Range("C15").Value = "Il6st"
$ws.Range("D15").Value = "MuSCs"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.4378713333333333
$ws.Range("H15").Value = 1.313614
$ws.Range("I15").Value = 0.3865237187563814
$ws.Range("J15").Value = 0.3865237187563814
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 13.37928666666667
$ws.Range("N15").Value = 40.13786
$ws.Range("O15").Value = 0.04830841776694433
$ws.Range("P15").Value = 0.04830841776694433
$ws.Range("Q15").Value = 5.858406091782221
$ws.Range("R15").Value = 52.72565482604
$ws.Range("S15").Value = 0.01867234928251617
$ws.Range("T15").Value = 0.01867234928251616

# Row 16
$ws.Range("A16").Value = "MuSCs"
$ws.Range("B16").Value = "Ctf1"
$ws.Range("C16").Value = "Il6st"
$ws.Range("D16").Value = "Resolving-Mac"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.4378713333333333
$ws.Range("H16").Value = 1.313614
$ws.Range("I16").Value = 0.3865237187563814
$ws.Range("J16").Value = 0.3865237187563814
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 10.83857433333333
$ws.Range("N16").Value = 32.515723
$ws.Range("O16").Value = 0.03913470052160829
$ws.Range("P16").Value = 0.03913470052160829
$ws.Range("Q16").Value = 4.74590099476911
$ws.Range("R16").Value = 42.713108952922
$ws.Range("S16").Value = 0.01512648997802934
$ws.Range("T16").Value = 0.01512648997802933

# Row 17
$ws.Range("A17").Value = "Resolving-Mac"
$ws.Range("B17").Value = "Ctf1"
$ws.Range("C17").Value = "Il6st"
$ws.Range("D17").Value = "ECs"
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 0.3333333333333333
$ws.Range("G17").Value = 0.013715
$ws.Range("H17").Value = 0.041145
$ws.Range("I17").Value = 0.01210669070840545
$ws.Range("J17").Value = 0.01210669070840545
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 108.845309
$ws.Range("N17").Value = 326.535927
$ws.Range("O17").Value = 0.3930063530400584
$ws.Range("P17").Value = 0.3930063530400583
$ws.Range("Q17").Value = 1.492813412935
$ws.Range("R17").Value = 13.435320716415
$ws.Range("S17").Value = 0.004758006362694386
$ws.Range("T17").Value = 0.004758006362694384

# Row 18
$ws.Range("A18").Value = "Resolving-Mac"
$ws.Range("B18").Value = "Ctf1"
$ws.Range("C18").Value = "Il6st"
$ws.Range("D18").Value = "FAPs"
$ws.Range("E18").Value = 1
$ws.Range("F18").Value = 0.3333333333333333
$ws.Range("G18").Value = 0.013715
$ws.Range("H18").Value = 0.041145
$ws.Range("I18").Value = 0.01210669070840545
$ws.Range("J18").Value = 0.01210669070840545
$ws.Range("K18").Value = 3
$ws.Range("L18").Value = 1
$ws.Range("M18").Value = 125.002688
$ws.Range("N18").Value = 375.008064
$ws.Range("O18").Value = 0.4513455929560021
$ws.Range("P18").Value = 0.451345592956002
$ws.Range("Q18").Value = 1.71441186592
$ws.Range("R18").Value = 15.42970679328
$ws.Range("S18").Value = 0.005464301496520178
$ws.Range("T18").Value = 0.005464301496520176

# Row 19
$ws.Range("A19").Value = "Resolving-Mac"
$ws.Range("B19").Value = "Ctf1"
$ws.Range("C19").Value = "Il6st"
$ws.Range("D19").Value = "Inflammatory-Mac"
$ws.Range("E19").Value = 1
$ws.Range("F19").Value = 0.3333333333333333
$ws.Range("G19").Value = 0.013715
$ws.Range("H19").Value = 0.041145
$ws.Range("I19").Value = 0.01210669070840545
$ws.Range("J19").Value = 0.01210669070840545
$ws.Range("K19").Value = 3
$ws.Range("L19").Value = 1
$ws.Range("M19").Value = 18.88973866666667
$ws.Range("N19").Value = 56.669216
$ws.Range("O19").Value = 0.06820493571538706
$ws.Range("P19").Value = 0.06820493571538705
$ws.Range("Q19").Value = 0.2590727658133333
$ws.Range("R19").Value = 2.33165489232
$ws.Range("S19").Value = 0.0008257360614928675
$ws.Range("T19").Value = 0.0008257360614928672

# Row 20
$ws.Range("A20").Value = "Resolving-Mac"
$ws.Range("B20").Value = "Ctf1"
$ws.Range("C20").Value = "Il6st"
$ws.Range("D20").Value = "MuSCs"
$ws.Range("E20").Value = 1
$ws.Range("F20").Value = 0.3333333333333333
$ws.Range("G20").Value = 0.013715
$ws.Range("H20").Value = 0.041145
$ws.Range("I20").Value = 0.01210669070840545
$ws.Range("J20").Value = 0.01210669070840545
$ws.Range("K20").Value = 3
$ws.Range("L20").Value = 1
$ws.Range("M20").Value = 13.37928666666667
$ws.Range("N20").Value = 40.13786
$ws.Range("O20").Value = 0.04830841776694433
$ws.Range("P20").Value = 0.04830841776694433
$ws.Range("Q20").Value = 0.1834969166333333
$ws.Range("R20").Value = 1.6514722497
$ws.Range("S20").Value = 0.0005848550725168337
$ws.Range("T20").Value = 0.0005848550725168335

# Row 21
$ws.Range("A21").Value = "Resolving-Mac"
$ws.Range("B21").Value = "Ctf1"
$ws.Range("C21").Value = "Il6st"
$ws.Range("D21").Value = "Resolving-Mac"
$ws.Range("E21").Value = 1
$ws.Range("F21").Value = 0.3333333333333333
$ws.Range("G21").Value = 0.013715
$ws.Range("H21").Value = 0.041145
$ws.Range("I21").Value = 0.01210669070840545
$ws.Range("J21").Value = 0.01210669070840545
$ws.Range("K21").Value = 3
$ws.Range("L21").Value = 1
$ws.Range("M21").Value = 10.83857433333333
$ws.Range("N21").Value = 32.515723
$ws.Range("O21").Value = 0.03913470052160829
$ws.Range("P21").Value = 0.03913470052160829
$ws.Range("Q21").Value = 0.1486510469816667
$ws.Range("R21").Value = 1.337859422835
$ws.Range("S21").Value = 0.000473791715181185
$ws.Range("T21").Value = 0.0004737917151811848

